$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9823903441429138
$ws.Range("B1").Value = 3.106423616409302
$ws.Range("C1").Value = 6.754115104675293
$ws.Range("D1").Value = 1.927229285240173
$ws.Range("E1").Value = 1.349624872207642
